# Weekly fruit/vegetable price update: insert a new price record for
# "Vega Modelo de Temuco - Brócoli" as row 279, pushing the existing
# rows 279-292 down to 280-293.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 279 (shifts old rows 279..292 -> 280..293)
$ws.Rows.Item(279).Insert()

# Match the date-number-format style used by the other rows' date column (D)
# before writing the value, so the cell keeps the same numeric date format.
$ws.Range("D279").NumberFormat = $ws.Range("D280").NumberFormat

# Populate the new row 279 with the new weekly record.
$ws.Range("A279").Value = 10
$ws.Range("B279").Value = "Vega Modelo de Temuco"
$ws.Range("C279").Value = "La Araucanía"
$ws.Range("D279").Value = 44509
$ws.Range("E279").Value = 9
$ws.Range("F279").Value = 100112023
$ws.Range("G279").Value = "Brócoli"
$ws.Range("H279").Value = "Sin especificar"
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 1100
$ws.Range("K279").Value = 800
$ws.Range("L279").Value = 900
$ws.Range("M279").Value = 845
$ws.Range("N279").Value = "$/unidad"
$ws.Range("O279").Value = "Región del Maule"
$ws.Range("P279").Value = 845
$ws.Range("Q279").Value = 1
$ws.Range("R279").Value = "Hortaliza"
